$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 379, shifting existing rows 379-394 down to 380-395.
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with its data.
$ws.Cells.Item(379, 1).Value = 10
$ws.Cells.Item(379, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(379, 3).Value = "La Araucanía"
$ws.Cells.Item(379, 4).Value2 = 45075
$ws.Cells.Item(379, 5).Value = 9
$ws.Cells.Item(379, 6).Value = 100112052
$ws.Cells.Item(379, 7).Value = "Albahaca"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 10).Value = 65
$ws.Cells.Item(379, 11).Value = 6000
$ws.Cells.Item(379, 12).Value = 6000
$ws.Cells.Item(379, 13).Value = 6000
$ws.Cells.Item(379, 14).Value = "$/paquete"
$ws.Cells.Item(379, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(379, 16).Value = 6000
$ws.Cells.Item(379, 17).Value = 1
$ws.Cells.Item(379, 18).Value = "Hortaliza"
